# GitNote.docx edit: "Add files via upload / add5"
#
# 1. Collapse the three runs (split by a pair of <w:proofErr/> grammar
#    markers) that make up the "test add remote" list item into a single
#    run, dropping the proofErr markers in the process.
# 2. Append a brand-new list item "add5" after the last paragraph
#    ("add remote4"), written as two runs ("a" / "dd5") matching the
#    existing "a" + "dd remote4" item's run split, and move the document's
#    "_GoBack" bookmark from the old last paragraph onto the new one.

$d = $word.ActiveDocument

# --- Change 1: merge "test add remote" runs ---------------------------
$d.Content.Find.Execute("test add remote", $true, $false, $false, $false, `
    $false, $true, 1, $false, "test add remote", 2) | Out-Null

# --- Change 2: add the "add5" list paragraph at the end ----------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# Range spanning from the start of the last paragraph through the very
# end of the document (this also covers the trailing _GoBack bookmark
# markers sitting in that paragraph), so the replacement XML below both
# rewrites the old paragraph (now bookmark-free) and introduces the new
# one carrying the bookmark.
$tailRange = $d.Range($lastPara.Range.Start, $d.Content.End)

$packageXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="a6"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
<w:ind w:firstLineChars="0"/>
</w:pPr>
<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>a</w:t></w:r>
<w:r><w:t>dd remote4</w:t></w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="a6"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
<w:ind w:firstLineChars="0"/>
</w:pPr>
<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>a</w:t></w:r>
<w:r><w:t>dd5</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$tailRange.InsertXML($packageXml) | Out-Null
